$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.512.58'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.804.48'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.94%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.010'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.65%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.009'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '308.20'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4547'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.71%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.61%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07123'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8710'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07783'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.19'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.837.02'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.274'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.322'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.87%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '86.29'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -6.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.011'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008568'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -4.54%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.559.94'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.12%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.84%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.956'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.054.84'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.35'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.985'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.18'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.71%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.85'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.999'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '112.78'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.863'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08690'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.086'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.70%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7304'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -5.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.433'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.85%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.41%  '
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.007'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.42%  '
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.506'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -7.44%  '
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.075'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.80%  '
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01914'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.39%  '
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05086'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.99%  '
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.860'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.26%  '
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.870'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.83%  '
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4896'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -4.87%  '
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1566'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.61%  '
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.122'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.47%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.009'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.56%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4587'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.74%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '101.67'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.68%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.887'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -4.84%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.581'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -4.39%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06000'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.59%  '
